$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Listado"

$ws.Range("A1").Value = "Mi listadito"
$ws.Range("B1").Value = "soy la celta b1"

$ws.Range("B2").Select()
